# Updated cryptos list on Fri Dec 29 09:49:51 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.743.08'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '2.359.07'
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.63%  '
$ws.Range("E7").Value = '  -2.27%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("E11").Value = '  -1.39%  '
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.96%  '
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.34%  '
$ws.Range("D16").Value = '2.718.03'
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("D17").Value = '2.299.81'
$ws.Range("E17").Value = '  -3.96%  '
$ws.Range("D18").Value = '42.728.30'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.78%  '
$ws.Range("E24").Value = '  -3.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.21%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("E29").Value = '  +1.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("E32").Value = '  -3.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.44%  '
$ws.Range("E35").Value = '  +16.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.131'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.70'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.98%  '
$ws.Range("E38").Value = '  -0.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.62%  '
$ws.Range("E41").Value = '  +2.57%  '
$ws.Range("E42").Value = '  -5.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '112.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.90%  '
$ws.Range("E47").Value = '  -0.11%  '

# Row 48 / 49: coin order swapped (FraxShare <-> BitcoinSV) with refreshed data
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.14%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.37%  '
$ws.Range("E51").Value = '  -2.06%  '